# Update existing data rows (2-16) with recalculated TPM values, and add new
# "Resolving-Mac" sending-cluster rows (17-21), per the updated NATMI output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Lama1", "Rpsa", "ECs", 2, 0.6666666666666666, 0.06628666666666667, 0.19886, 0.09845363529874238, 0.09845363529874238, 3, 1, 127.5808283333333, 382.742485, 0.1532286472569342, 0.1532286472569342, 8.45690784078889, 76.11217056709999, 0.01508591735435384, 0.01508591735435384),
    @("ECs", "Lama1", "Rpsa", "FAPs", 2, 0.6666666666666666, 0.06628666666666667, 0.19886, 0.09845363529874238, 0.09845363529874238, 3, 1, 284.3623303333333, 853.0869909999999, 0.341528235684153, 0.341528235684153, 18.84943100336222, 169.64487903026, 0.03362469636027053, 0.03362469636027053),
    @("ECs", "Lama1", "Rpsa", "Inflammatory-Mac", 2, 0.6666666666666666, 0.06628666666666667, 0.19886, 0.09845363529874238, 0.09845363529874238, 3, 1, 155.9545746666667, 467.863724, 0.1873064223040503, 0.1873064223040503, 10.33770890607111, 93.03938015464001, 0.01844099819063519, 0.01844099819063519),
    @("ECs", "Lama1", "Rpsa", "MuSCs", 2, 0.6666666666666666, 0.06628666666666667, 0.19886, 0.09845363529874238, 0.09845363529874238, 3, 1, 132.7811556666666, 398.343467, 0.1594744063806953, 0.1594744063806954, 8.80162020529111, 79.21458184762, 0.01570083504528842, 0.01570083504528842),
    @("ECs", "Lama1", "Rpsa", "Resolving-Mac", 2, 0.6666666666666666, 0.06628666666666667, 0.19886, 0.09845363529874238, 0.09845363529874238, 3, 1, 131.9384486666667, 395.815346, 0.1584622883741672, 0.1584622883741672, 8.745759967284444, 78.71183970556, 0.0156011883481944, 0.0156011883481944),
    @("FAPs", "Lama1", "Rpsa", "ECs", 3, 1, 0.585754, 1.757262, 0.8700031784790234, 0.8700031784790236, 3, 1, 127.5808283333333, 382.742485, 0.1532286472569342, 0.1532286472569342, 74.73098051956333, 672.57882467607, 0.1333094101475738, 0.1333094101475739),
    @("FAPs", "Lama1", "Rpsa", "FAPs", 3, 1, 0.585754, 1.757262, 0.8700031784790234, 0.8700031784790236, 3, 1, 284.3623303333333, 853.0869909999999, 0.341528235684153, 0.341528235684153, 166.5663724420713, 1499.097351978642, 0.2971306505855461, 0.2971306505855462),
    @("FAPs", "Lama1", "Rpsa", "Inflammatory-Mac", 3, 1, 0.585754, 1.757262, 0.8700031784790234, 0.8700031784790236, 3, 1, 155.9545746666667, 467.863724, 0.1873064223040503, 0.1873064223040503, 91.35101592929867, 822.1591433636881, 0.162957182754058, 0.162957182754058),
    @("FAPs", "Lama1", "Rpsa", "MuSCs", 3, 1, 0.585754, 1.757262, 0.8700031784790234, 0.8700031784790236, 3, 1, 132.7811556666666, 398.343467, 0.1594744063806953, 0.1594744063806954, 77.77709305637265, 699.993837507354, 0.1387432404372604, 0.1387432404372605),
    @("FAPs", "Lama1", "Rpsa", "Resolving-Mac", 3, 1, 0.585754, 1.757262, 0.8700031784790234, 0.8700031784790236, 3, 1, 131.9384486666667, 395.815346, 0.1584622883741672, 0.1584622883741672, 77.28347406029467, 695.551266542652, 0.1378626945545851, 0.1378626945545851),
    @("MuSCs", "Lama1", "Rpsa", "ECs", 2, 0.6666666666666666, 0.01811333333333333, 0.05434, 0.02690320095611817, 0.02690320095611818, 3, 1, 127.5808283333333, 382.742485, 0.1532286472569342, 0.1532286472569342, 2.310914070544444, 20.7982266349, 0.004122341089387446, 0.004122341089387447),
    @("MuSCs", "Lama1", "Rpsa", "FAPs", 2, 0.6666666666666666, 0.01811333333333333, 0.05434, 0.02690320095611817, 0.02690320095611818, 3, 1, 284.3623303333333, 853.0869909999999, 0.341528235684153, 0.341528235684153, 5.150749676771111, 46.35674709094, 0.009188202756799256, 0.009188202756799258),
    @("MuSCs", "Lama1", "Rpsa", "Inflammatory-Mac", 2, 0.6666666666666666, 0.01811333333333333, 0.05434, 0.02690320095611817, 0.02690320095611818, 3, 1, 155.9545746666667, 467.863724, 0.1873064223040503, 0.1873064223040503, 2.824857195795555, 25.42371476216, 0.005039142319617399, 0.0050391423196174),
    @("MuSCs", "Lama1", "Rpsa", "MuSCs", 2, 0.6666666666666666, 0.01811333333333333, 0.05434, 0.02690320095611817, 0.02690320095611818, 3, 1, 132.7811556666666, 398.343467, 0.1594744063806953, 0.1594744063806954, 2.405109332975555, 21.64598399678, 0.004290372002217501, 0.004290372002217503),
    @("MuSCs", "Lama1", "Rpsa", "Resolving-Mac", 2, 0.6666666666666666, 0.01811333333333333, 0.05434, 0.02690320095611817, 0.02690320095611818, 3, 1, 131.9384486666667, 395.815346, 0.1584622883741672, 0.1584622883741672, 2.389845100182222, 21.50860590164, 0.004263142788096568, 0.004263142788096569),
    @("Resolving-Mac", "Lama1", "Rpsa", "ECs", 1, 0.3333333333333333, 0.003124, 0.009372, 0.004639985266115928, 0.004639985266115929, 3, 1, 127.5808283333333, 382.742485, 0.1532286472569342, 0.1532286472569342, 0.3985625077133333, 3.58706256942, 0.0007109786656190495, 0.0007109786656190496),
    @("Resolving-Mac", "Lama1", "Rpsa", "FAPs", 1, 0.3333333333333333, 0.003124, 0.009372, 0.004639985266115928, 0.004639985266115929, 3, 1, 284.3623303333333, 853.0869909999999, 0.341528235684153, 0.341528235684153, 0.8883479199613333, 7.995131279651999, 0.001584685981537038, 0.001584685981537038),
    @("Resolving-Mac", "Lama1", "Rpsa", "Inflammatory-Mac", 1, 0.3333333333333333, 0.003124, 0.009372, 0.004639985266115928, 0.004639985266115929, 3, 1, 155.9545746666667, 467.863724, 0.1873064223040503, 0.1873064223040503, 0.4872020912586667, 4.384818821328, 0.0008690990397396811, 0.0008690990397396813),
    @("Resolving-Mac", "Lama1", "Rpsa", "MuSCs", 1, 0.3333333333333333, 0.003124, 0.009372, 0.004639985266115928, 0.004639985266115929, 3, 1, 132.7811556666666, 398.343467, 0.1594744063806953, 0.1594744063806954, 0.4148083303026666, 3.733274972724, 0.0007399588959290103, 0.0007399588959290106),
    @("Resolving-Mac", "Lama1", "Rpsa", "Resolving-Mac", 1, 0.3333333333333333, 0.003124, 0.009372, 0.004639985266115928, 0.004639985266115929, 3, 1, 131.9384486666667, 395.815346, 0.1584622883741672, 0.1584622883741672, 0.4121757136346667, 3.709581422712, 0.000735262683291149, 0.0007352626832911492)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}